# Update the "Förändrad" (Changed) date column C for rows 2-261 from
# serial date 45180 to 45181 (i.e. bump the date by one day), matching
# the upstream source update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C261").Value = 45181
